# Generate Report for Handback
#
# For each locale sheet (zh-cn, de-de):
#  - Status column (C) moves from "Ready for handoff" to
#    "Handed back: in sync with en-US"
#  - A "Latest Target File" (F) and "Latest Handback File" (G) hyperlinked
#    file name are recorded - they mirror the existing Source (A) / Latest
#    Handoff File (D) links, since the handback target/content matches the
#    handoff.
#  - The "Latest Handback DateTime" (H) placeholder ("0001-01-01 00:00:00")
#    is replaced with the real handback timestamp for that locale.

$wb = $excel.ActiveWorkbook

$statusText = "Handed back: in sync with en-US"

$locales = @(
    @{ Sheet = "zh-cn"; HandbackTime = "2016-03-18 10:46:33" },
    @{ Sheet = "de-de"; HandbackTime = "2016-03-18 10:46:38" }
)

foreach ($locale in $locales) {
    $ws = $wb.Worksheets.Item($locale.Sheet)

    # Collect the existing hyperlink addresses/text for row 2 and row 3
    # before we start mutating the sheet, keyed by the anchor address.
    $linkInfo = @{}
    foreach ($hl in $ws.Hyperlinks) {
        $addr = $hl.Range.Address()
        $linkInfo[$addr] = @{ Address = $hl.Address; Text = $hl.TextToDisplay }
    }

    $aLinks = @{ 2 = $linkInfo["`$A`$2"]; 3 = $linkInfo["`$A`$3"] }
    $dLinks = @{ 2 = $linkInfo["`$D`$2"]; 3 = $linkInfo["`$D`$3"] }

    foreach ($row in 2, 3) {
        # Status text: "Ready for handoff" -> "Handed back: in sync with en-US"
        $ws.Range("C$row").Value = $statusText

        # F = Latest Target File : same file/link as column A (source file)
        $aInfo = $aLinks[$row]
        $ws.Hyperlinks.Add($ws.Range("F$row"), $aInfo.Address, "", "", $aInfo.Text) | Out-Null
        $ws.Range("F$row").Font.Underline = 2
        $ws.Range("F$row").Font.Color = 15570276
        $ws.Range("F$row").Font.Name = "Calibri"
        $ws.Range("F$row").Font.Size = 11

        # G = Latest Handback File : same file/link as column D (handoff file)
        $dInfo = $dLinks[$row]
        $ws.Hyperlinks.Add($ws.Range("G$row"), $dInfo.Address, "", "", $dInfo.Text) | Out-Null
        $ws.Range("G$row").Font.Underline = 2
        $ws.Range("G$row").Font.Color = 15570276
        $ws.Range("G$row").Font.Name = "Calibri"
        $ws.Range("G$row").Font.Size = 11

        # H = Latest Handback DateTime : replace the placeholder date
        $ws.Range("H$row").Value = $locale.HandbackTime
    }
}

Write-Host "Handback report generated"
